# Update odds data in Sheet1 for the 2024-11-13 FlashScore export.
# Row 2: Flamengo RJ vs Atletico-MG
# Row 6: Carabobo vs Monagas

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.91
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 4.75
$ws.Range("Z2").Value = 15
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 15
$ws.Range("AR2").Value = 51
$ws.Range("AW2").Value = 6
$ws.Range("BB2").Value = 301

# --- Row 6 updates ---
$ws.Range("H6").Value = 3.45
$ws.Range("I6").Value = 7.2
$ws.Range("J6").Value = 2.05
$ws.Range("K6").Value = 2.07
$ws.Range("L6").Value = 6.6
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.35
$ws.Range("P6").Value = 2.7
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.62
$ws.Range("S6").Value = 1.42
$ws.Range("T6").Value = 2.47
$ws.Range("U6").Value = 2.07
$ws.Range("V6").Value = 1.6
$ws.Range("W6").Value = 5.2
$ws.Range("X6").Value = 6.1
$ws.Range("Z6").Value = 10.5
$ws.Range("AA6").Value = 14
$ws.Range("AB6").Value = 35
$ws.Range("AC6").Value = 7.8
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 19.5
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 50
$ws.Range("AK6").Value = 200
$ws.Range("AO6").Value = 7.2
$ws.Range("AR6").Value = 60
$ws.Range("AT6").Value = 2.45
$ws.Range("AU6").Value = 7.8
$ws.Range("AW6").Value = 8
$ws.Range("AX6").Value = 45
